# Daily attendance processing - 2026-01-04 17:01:21
# Normalizes the "Recorded By" (column G) values so that "System" is
# listed before the recorder's email address, e.g.
#   "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Columns.Item(7)   # Column G - "Recorded By"
$col.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com")
